$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 09:06"

# Row 19: Argentina
$ws.Range("A19").Value = "Argentina"
$ws.Range("B19").Value = 253868
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 174974
$ws.Range("E19").Value = 74130
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 4764

# Row 29: Kazajistan
$ws.Range("A29").Value = "Kazajistan"
$ws.Range("B29").Value = 100164
$ws.Range("C29").Value = 722
$ws.Range("D29").Value = 73702
$ws.Range("E29").Value = 25193
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 23
$ws.Range("H29").Value = 1269

# Row 35: Ucrania
$ws.Range("A35").Value = "Ucrania"
$ws.Range("B35").Value = 83115
$ws.Range("C35").Value = 1158
$ws.Range("D35").Value = 44934
$ws.Range("E35").Value = 36230
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 29
$ws.Range("H35").Value = 1951

# Row 36: Suecia
$ws.Range("A36").Value = "Suecia"
$ws.Range("B36").Value = 82972
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 5766

# Row 55: Armenia
$ws.Range("A55").Value = "Armenia"
$ws.Range("B55").Value = 40593
$ws.Range("C55").Value = 160
$ws.Range("D55").Value = 33157
$ws.Range("E55").Value = 6633
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 7
$ws.Range("H55").Value = 803

# Row 56: Kirguistan
$ws.Range("A56").Value = "Kirguistan"
$ws.Range("B56").Value = 40455
$ws.Range("C56").Value = 278
$ws.Range("D56").Value = 32734
$ws.Range("E56").Value = 6243
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = 1478

# Row 94: Haiti
$ws.Range("A94").Value = "Haiti"
$ws.Range("B94").Value = 7649
$ws.Range("C94").Value = 15
$ws.Range("D94").Value = 4982
$ws.Range("E94").Value = 2484
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 183

# Row 142: Letonia
$ws.Range("A142").Value = "Letonia"
$ws.Range("B142").Value = 1293
$ws.Range("C142").Value = 3
$ws.Range("D142").Value = 1078
$ws.Range("E142").Value = 183
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 32

# Row 144: Georgia
$ws.Range("A144").Value = "Georgia"
$ws.Range("B144").Value = 1264
$ws.Range("C144").Value = 14
$ws.Range("D144").Value = 1054
$ws.Range("E144").Value = 193
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 17

# Row 145: Siria
$ws.Range("A145").Value = "Siria"
$ws.Range("B145").Value = 1255
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 364
$ws.Range("E145").Value = 839
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 52

# Row 146: Republica de Chipre
$ws.Range("A146").Value = "Republica de Chipre"
$ws.Range("B146").Value = 1252
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 870
$ws.Range("E146").Value = 363
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 19

# Row 175: Islas Feroe
$ws.Range("A175").Value = "Islas Feroe"
$ws.Range("B175").Value = 313
$ws.Range("C175").Value = 7
$ws.Range("D175").Value = 225
$ws.Range("E175").Value = 88
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

# Row 203: Santa Lucia
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("B203").Value = 25
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 25
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 213: Islas Malvinas
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214: Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

